# Update "想去人数" (want-to-go count) figures in the 广州-漫展信息 workbook
# to the values captured at the newer scrape (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 574
$ws1.Range("F3").Value  = 255
$ws1.Range("F4").Value  = 29
$ws1.Range("F5").Value  = 737
$ws1.Range("F6").Value  = 360
$ws1.Range("F8").Value  = 145
$ws1.Range("F10").Value = 213
$ws1.Range("F11").Value = 5910
$ws1.Range("F13").Value = 40
$ws1.Range("F14").Value = 490
$ws1.Range("F16").Value = 547
$ws1.Range("F17").Value = 354
$ws1.Range("F18").Value = 421
$ws1.Range("F21").Value = 704
$ws1.Range("F22").Value = 130
$ws1.Range("F24").Value = 307
$ws1.Range("F25").Value = 1016
$ws1.Range("F27").Value = 1794
$ws1.Range("F28").Value = 461
$ws1.Range("F29").Value = 32

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 51
$ws2.Range("F6").Value = 299

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 216

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 216
$ws4.Range("F3").Value  = 574
$ws4.Range("F4").Value  = 255
$ws4.Range("F5").Value  = 29
$ws4.Range("F6").Value  = 737
$ws4.Range("F8").Value  = 360
$ws4.Range("F10").Value = 145
$ws4.Range("F12").Value = 213
$ws4.Range("F13").Value = 5910
$ws4.Range("F15").Value = 40
$ws4.Range("F17").Value = 490
$ws4.Range("F19").Value = 547
$ws4.Range("F20").Value = 354
$ws4.Range("F21").Value = 421
$ws4.Range("F22").Value = 51
$ws4.Range("F26").Value = 299
$ws4.Range("F28").Value = 704
$ws4.Range("F32").Value = 130
$ws4.Range("F34").Value = 307
$ws4.Range("F35").Value = 1016
$ws4.Range("F37").Value = 1794
$ws4.Range("F38").Value = 461
$ws4.Range("F39").Value = 32

$wb.Save()
